$wb = $excel.ActiveWorkbook

# 1. Add a new "Supervisors" sheet as the first sheet in the workbook.
$supervisors = $wb.Worksheets.Add()
$supervisors.Name = "Supervisors"

$supervisors.Range("A1").Value = "Supervisor"
$supervisors.Range("B1").Value = "Max_number_of_projects"
$supervisors.Range("C1").Value = "Max_number_of_students"
$supervisors.Range("A2").Value = "Dr Smith"

# 2. Update the "Projects" sheet: rename header + add Supervisor column.
$projects = $wb.Worksheets.Item("Projects")
$projects.Range("A1").Value = "Project"
$projects.Range("C1").Value = "Supervisor"
$projects.Range("C2").Value = "Dr Smith"
$projects.Range("C3").Value = "Dr Smith"
$projects.Range("C4").Value = "Dr Smith"
$projects.Range("C5").Value = "Dr Smith"
$projects.Range("C6").Value = "Dr Smith"

# 3. Make the Supervisors sheet the active sheet / tab.
$wb.Worksheets.Item("Supervisors").Activate()
$supervisors.Range("A1").Select()
